# Add results analysis for user study:
#  - Fill in missing encoding-scheme answers for participant 1 (row 2, cols K:P)
#  - Remove participant 13's row (row 21) of responses, keeping the
#    timestamp cells' date formatting but blanking their values
#  - Update the sheet selection to A24 (scrolled back to top of sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank "Agree" style answers for participant 1
$ws.Range("K2").Value2 = "Agree"
$ws.Range("L2").Value2 = "Strongly Agree"
$ws.Range("M2").Value2 = "Strongly Disagree"
$ws.Range("N2").Value2 = "Agree"
$ws.Range("O2").Value2 = "Agree"
$ws.Range("P2").Value2 = "Strongly Agree"

# Clear out participant 13's row of survey answers entirely, leaving the
# (date-formatted) start/completion time cells blank but still formatted
$ws.Range("A21:AZ21").ClearContents()

# Reset the view: scroll back to the top and select A24
$ws.Range("A24").Select()
